$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Correct two pairs of rows whose match data was swapped (74/75 and 77/78) ---
# For each pair, columns F:V (match details) are exchanged while A:E (index/meta) stay put.
# Row 74 gets the match info previously stored in row 75 (Mafra vs Leixoes)
$ws.Cells.Item(74, 6).Value = "Mafra"
$ws.Cells.Item(74, 7).Value = 0
$ws.Cells.Item(74, 8).Value = "Leixoes"
$ws.Cells.Item(74, 9).Value = 1
$ws.Cells.Item(74, 10).Value = 1.88
$ws.Cells.Item(74, 11).Value = "01/11/2023 16:12"
$ws.Cells.Item(74, 12).Value = 1.93
$ws.Cells.Item(74, 13).Value = "04/11/2023 11:48"
$ws.Cells.Item(74, 14).Value = 3.73
$ws.Cells.Item(74, 15).Value = "01/11/2023 16:12"
$ws.Cells.Item(74, 16).Value = 3.55
$ws.Cells.Item(74, 17).Value = "04/11/2023 11:51"
$ws.Cells.Item(74, 18).Value = 3.8
$ws.Cells.Item(74, 19).Value = "01/11/2023 16:12"
$ws.Cells.Item(74, 20).Value = 4.21
$ws.Cells.Item(74, 21).Value = "04/11/2023 11:51"
$ws.Cells.Item(74, 22).Value = "https://www.betexplorer.com/football/portugal/liga-portugal-2/mafra-leixoes/YiBBPnTT/"

# Row 75 gets the match info previously stored in row 74 (FC Porto B vs Feirense)
$ws.Cells.Item(75, 6).Value = "FC Porto B"
$ws.Cells.Item(75, 7).Value = 2
$ws.Cells.Item(75, 8).Value = "Feirense"
$ws.Cells.Item(75, 9).Value = 0
$ws.Cells.Item(75, 10).Value = 1.98
$ws.Cells.Item(75, 11).Value = "01/11/2023 16:12"
$ws.Cells.Item(75, 12).Value = 1.84
$ws.Cells.Item(75, 13).Value = "04/11/2023 11:59"
$ws.Cells.Item(75, 14).Value = 3.57
$ws.Cells.Item(75, 15).Value = "01/11/2023 16:12"
$ws.Cells.Item(75, 16).Value = 3.78
$ws.Cells.Item(75, 17).Value = "04/11/2023 11:59"
$ws.Cells.Item(75, 18).Value = 3.87
$ws.Cells.Item(75, 19).Value = "01/11/2023 16:12"
$ws.Cells.Item(75, 20).Value = 4.41
$ws.Cells.Item(75, 21).Value = "04/11/2023 11:58"
$ws.Cells.Item(75, 22).Value = "https://www.betexplorer.com/football/portugal/liga-portugal-2/fc-porto-feirense/jTL6QSDN/"

# Row 77 gets the match info previously stored in row 78 (Nacional vs Santa Clara)
$ws.Cells.Item(77, 6).Value = "Nacional"
$ws.Cells.Item(77, 7).Value = 1
$ws.Cells.Item(77, 8).Value = "Santa Clara"
$ws.Cells.Item(77, 9).Value = 1
$ws.Cells.Item(77, 10).Value = 2.98
$ws.Cells.Item(77, 11).Value = "01/11/2023 16:12"
$ws.Cells.Item(77, 12).Value = 2.81
$ws.Cells.Item(77, 13).Value = "04/11/2023 18:58"
$ws.Cells.Item(77, 14).Value = 3.27
$ws.Cells.Item(77, 15).Value = "01/11/2023 16:12"
$ws.Cells.Item(77, 16).Value = 3.23
$ws.Cells.Item(77, 17).Value = "04/11/2023 18:52"
$ws.Cells.Item(77, 18).Value = 2.39
$ws.Cells.Item(77, 19).Value = "01/11/2023 16:12"
$ws.Cells.Item(77, 20).Value = 2.72
$ws.Cells.Item(77, 21).Value = "04/11/2023 18:52"
$ws.Cells.Item(77, 22).Value = "https://www.betexplorer.com/football/portugal/liga-portugal-2/nacional-santa-clara/xQH2R8bH/"

# Row 78 gets the match info previously stored in row 77 (Benfica B vs Penafiel)
$ws.Cells.Item(78, 6).Value = "Benfica B"
$ws.Cells.Item(78, 7).Value = 1
$ws.Cells.Item(78, 8).Value = "Penafiel"
$ws.Cells.Item(78, 9).Value = 0
$ws.Cells.Item(78, 10).Value = 1.91
$ws.Cells.Item(78, 11).Value = "29/10/2023 16:42"
$ws.Cells.Item(78, 12).Value = 2.26
$ws.Cells.Item(78, 13).Value = "04/11/2023 18:53"
$ws.Cells.Item(78, 14).Value = 3.67
$ws.Cells.Item(78, 15).Value = "29/10/2023 16:42"
$ws.Cells.Item(78, 16).Value = 3.59
$ws.Cells.Item(78, 17).Value = "04/11/2023 18:53"
$ws.Cells.Item(78, 18).Value = 4.01
$ws.Cells.Item(78, 19).Value = "29/10/2023 16:42"
$ws.Cells.Item(78, 20).Value = 3.2
$ws.Cells.Item(78, 21).Value = "04/11/2023 18:53"
$ws.Cells.Item(78, 22).Value = "https://www.betexplorer.com/football/portugal/liga-portugal-2/benfica-penafiel/xjmbUAEb/"

# --- Append 7 new match rows (92-98) ---
# Seed formatting for the new rows by copying the style of the last existing row (91),
# then overwrite the values cell-by-cell so the shared-style indices (bold index column,
# date-formatted E column) match the rest of the table without touching styles.xml.
$ws.Range("A91:V91").Copy() | Out-Null
$ws.Range("A92:V98").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Row 92
$ws.Cells.Item(92, 1).Value = 91
$ws.Cells.Item(92, 2).Value = "portugal"
$ws.Cells.Item(92, 3).Value = "liga-portugal-2"
$ws.Cells.Item(92, 4).Value = "2023-2024"
$ws.Cells.Item(92, 5).Value = 45248.5
$ws.Cells.Item(92, 6).Value = "Mafra"
$ws.Cells.Item(92, 7).Value = 1
$ws.Cells.Item(92, 8).Value = "Vilaverdense"
$ws.Cells.Item(92, 9).Value = 0
$ws.Cells.Item(92, 10).Value = 1.63
$ws.Cells.Item(92, 11).Value = "14/11/2023 13:12"
$ws.Cells.Item(92, 12).Value = 1.72
$ws.Cells.Item(92, 13).Value = "18/11/2023 11:50"
$ws.Cells.Item(92, 14).Value = 3.96
$ws.Cells.Item(92, 15).Value = "14/11/2023 13:12"
$ws.Cells.Item(92, 16).Value = 3.89
$ws.Cells.Item(92, 17).Value = "18/11/2023 11:50"
$ws.Cells.Item(92, 18).Value = 5.48
$ws.Cells.Item(92, 19).Value = "14/11/2023 13:12"
$ws.Cells.Item(92, 20).Value = 5.03
$ws.Cells.Item(92, 21).Value = "18/11/2023 11:50"
$ws.Cells.Item(92, 22).Value = "https://www.betexplorer.com/football/portugal/liga-portugal-2/mafra-vilaverdense-fc/bZ43dTci/"

# Row 93
$ws.Cells.Item(93, 1).Value = 92
$ws.Cells.Item(93, 2).Value = "portugal"
$ws.Cells.Item(93, 3).Value = "liga-portugal-2"
$ws.Cells.Item(93, 4).Value = "2023-2024"
$ws.Cells.Item(93, 5).Value = 45248.625
$ws.Cells.Item(93, 6).Value = "Pacos Ferreira"
$ws.Cells.Item(93, 7).Value = 0
$ws.Cells.Item(93, 8).Value = "Santa Clara"
$ws.Cells.Item(93, 9).Value = 2
$ws.Cells.Item(93, 10).Value = 2.35
$ws.Cells.Item(93, 11).Value = "14/11/2023 13:12"
$ws.Cells.Item(93, 12).Value = 2.69
$ws.Cells.Item(93, 13).Value = "18/11/2023 14:59"
$ws.Cells.Item(93, 14).Value = 3.22
$ws.Cells.Item(93, 15).Value = "14/11/2023 13:12"
$ws.Cells.Item(93, 16).Value = 3.06
$ws.Cells.Item(93, 17).Value = "18/11/2023 14:53"
$ws.Cells.Item(93, 18).Value = 3.26
$ws.Cells.Item(93, 19).Value = "14/11/2023 13:12"
$ws.Cells.Item(93, 20).Value = 2.99
$ws.Cells.Item(93, 21).Value = "18/11/2023 14:59"
$ws.Cells.Item(93, 22).Value = "https://www.betexplorer.com/football/portugal/liga-portugal-2/pacos-ferreira-santa-clara/IsKmyJRr/"

# Row 94
$ws.Cells.Item(94, 1).Value = 93
$ws.Cells.Item(94, 2).Value = "portugal"
$ws.Cells.Item(94, 3).Value = "liga-portugal-2"
$ws.Cells.Item(94, 4).Value = "2023-2024"
$ws.Cells.Item(94, 5).Value = 45248.6875
$ws.Cells.Item(94, 6).Value = "Tondela"
$ws.Cells.Item(94, 7).Value = 2
$ws.Cells.Item(94, 8).Value = "Feirense"
$ws.Cells.Item(94, 9).Value = 0
$ws.Cells.Item(94, 10).Value = 1.73
$ws.Cells.Item(94, 11).Value = "14/11/2023 13:12"
$ws.Cells.Item(94, 12).Value = 1.94
$ws.Cells.Item(94, 13).Value = "18/11/2023 16:20"
$ws.Cells.Item(94, 14).Value = 3.74
$ws.Cells.Item(94, 15).Value = "14/11/2023 13:12"
$ws.Cells.Item(94, 16).Value = 3.45
$ws.Cells.Item(94, 17).Value = "18/11/2023 16:20"
$ws.Cells.Item(94, 18).Value = 4.54
$ws.Cells.Item(94, 19).Value = "14/11/2023 13:12"
$ws.Cells.Item(94, 20).Value = 4.34
$ws.Cells.Item(94, 21).Value = "18/11/2023 16:20"
$ws.Cells.Item(94, 22).Value = "https://www.betexplorer.com/football/portugal/liga-portugal-2/tondela-feirense/C4Jizwsk/"

# Row 95
$ws.Cells.Item(95, 1).Value = 94
$ws.Cells.Item(95, 2).Value = "portugal"
$ws.Cells.Item(95, 3).Value = "liga-portugal-2"
$ws.Cells.Item(95, 4).Value = "2023-2024"
$ws.Cells.Item(95, 5).Value = 45248.79166666666
$ws.Cells.Item(95, 6).Value = "Leiria"
$ws.Cells.Item(95, 7).Value = 0
$ws.Cells.Item(95, 8).Value = "Leixoes"
$ws.Cells.Item(95, 9).Value = 0
$ws.Cells.Item(95, 10).Value = 1.59
$ws.Cells.Item(95, 11).Value = "14/11/2023 13:12"
$ws.Cells.Item(95, 12).Value = 1.83
$ws.Cells.Item(95, 13).Value = "18/11/2023 18:50"
$ws.Cells.Item(95, 14).Value = 4.13
$ws.Cells.Item(95, 15).Value = "14/11/2023 13:12"
$ws.Cells.Item(95, 16).Value = 3.78
$ws.Cells.Item(95, 17).Value = "18/11/2023 18:50"
$ws.Cells.Item(95, 18).Value = 5.68
$ws.Cells.Item(95, 19).Value = "14/11/2023 13:12"
$ws.Cells.Item(95, 20).Value = 4.44
$ws.Cells.Item(95, 21).Value = "18/11/2023 18:50"
$ws.Cells.Item(95, 22).Value = "https://www.betexplorer.com/football/portugal/liga-portugal-2/leiria-leixoes/Ot2FgRsA/"

# Row 96
$ws.Cells.Item(96, 1).Value = 95
$ws.Cells.Item(96, 2).Value = "portugal"
$ws.Cells.Item(96, 3).Value = "liga-portugal-2"
$ws.Cells.Item(96, 4).Value = "2023-2024"
$ws.Cells.Item(96, 5).Value = 45249.5
$ws.Cells.Item(96, 6).Value = "AVS"
$ws.Cells.Item(96, 7).Value = 0
$ws.Cells.Item(96, 8).Value = "Nacional"
$ws.Cells.Item(96, 9).Value = 1
$ws.Cells.Item(96, 10).Value = 1.92
$ws.Cells.Item(96, 11).Value = "14/11/2023 13:12"
$ws.Cells.Item(96, 12).Value = 2.18
$ws.Cells.Item(96, 13).Value = "19/11/2023 11:58"
$ws.Cells.Item(96, 14).Value = 3.38
$ws.Cells.Item(96, 15).Value = "14/11/2023 13:12"
$ws.Cells.Item(96, 16).Value = 3.56
$ws.Cells.Item(96, 17).Value = "19/11/2023 11:58"
$ws.Cells.Item(96, 18).Value = 4.23
$ws.Cells.Item(96, 19).Value = "14/11/2023 13:12"
$ws.Cells.Item(96, 20).Value = 3.32
$ws.Cells.Item(96, 21).Value = "19/11/2023 11:51"
$ws.Cells.Item(96, 22).Value = "https://www.betexplorer.com/football/portugal/liga-portugal-2/avs-nacional/xQ37emDc/"

# Row 97
$ws.Cells.Item(97, 1).Value = 96
$ws.Cells.Item(97, 2).Value = "portugal"
$ws.Cells.Item(97, 3).Value = "liga-portugal-2"
$ws.Cells.Item(97, 4).Value = "2023-2024"
$ws.Cells.Item(97, 5).Value = 45249.625
$ws.Cells.Item(97, 6).Value = "Os Belenenses"
$ws.Cells.Item(97, 7).Value = 2
$ws.Cells.Item(97, 8).Value = "Penafiel"
$ws.Cells.Item(97, 9).Value = 0
$ws.Cells.Item(97, 10).Value = 2.49
$ws.Cells.Item(97, 11).Value = "12/11/2023 15:12"
$ws.Cells.Item(97, 12).Value = 2.77
$ws.Cells.Item(97, 13).Value = "19/11/2023 14:54"
$ws.Cells.Item(97, 14).Value = 3.19
$ws.Cells.Item(97, 15).Value = "12/11/2023 15:12"
$ws.Cells.Item(97, 16).Value = 3.11
$ws.Cells.Item(97, 17).Value = "19/11/2023 14:51"
$ws.Cells.Item(97, 18).Value = 3.05
$ws.Cells.Item(97, 19).Value = "12/11/2023 15:12"
$ws.Cells.Item(97, 20).Value = 2.86
$ws.Cells.Item(97, 21).Value = "19/11/2023 14:54"
$ws.Cells.Item(97, 22).Value = "https://www.betexplorer.com/football/portugal/liga-portugal-2/cf-os-belenenses-penafiel/65iKhocG/"

# Row 98
$ws.Cells.Item(98, 1).Value = 97
$ws.Cells.Item(98, 2).Value = "portugal"
$ws.Cells.Item(98, 3).Value = "liga-portugal-2"
$ws.Cells.Item(98, 4).Value = "2023-2024"
$ws.Cells.Item(98, 5).Value = 45249.6875
$ws.Cells.Item(98, 6).Value = "Maritimo"
$ws.Cells.Item(98, 7).Value = 1
$ws.Cells.Item(98, 8).Value = "Torreense"
$ws.Cells.Item(98, 9).Value = 2
$ws.Cells.Item(98, 10).Value = 1.89
$ws.Cells.Item(98, 11).Value = "13/11/2023 15:12"
$ws.Cells.Item(98, 12).Value = 1.69
$ws.Cells.Item(98, 13).Value = "19/11/2023 16:23"
$ws.Cells.Item(98, 14).Value = 3.47
$ws.Cells.Item(98, 15).Value = "13/11/2023 15:12"
$ws.Cells.Item(98, 16).Value = 3.79
$ws.Cells.Item(98, 17).Value = "19/11/2023 16:23"
$ws.Cells.Item(98, 18).Value = 4.05
$ws.Cells.Item(98, 19).Value = "13/11/2023 15:12"
$ws.Cells.Item(98, 20).Value = 5.47
$ws.Cells.Item(98, 21).Value = "19/11/2023 16:23"
$ws.Cells.Item(98, 22).Value = "https://www.betexplorer.com/football/portugal/liga-portugal-2/maritimo-torreense/rRMGYP5p/"

